# Swap the contents of column A (分类名称/classification name) and
# column B (单品名称/item name) for the header row and all data rows.
# This reorders the columns so that the item name now comes first
# (column A) and the classification name comes second (column B),
# matching the "前后顺序" (front/back order) change described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $valA = $ws.Cells.Item($r, 1).Value2
    $valB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $valB
    $ws.Cells.Item($r, 2).Value = $valA
}
